$wb = $excel.ActiveWorkbook

# Update the "想去人数" (F column) counts on both the "展览" sheet and the
# "全部类型" sheet, which mirror the same rows of data.
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 221
    $ws.Range("F3").Value = 235
}
